# CK_OBS_WITH_FEATURES.xlsx edit
# - Rename Sheet1 -> Ck_Old_Model
# - Move the active selection from F4 to E2
# - Widen column D (4th column) from 13 to 15 characters

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Ck_Old_Model"

# Move active selection to E2
$ws.Range("E2").Select()

# Resize column D so the stored OOXML width becomes 15
# (this engine's ColumnWidth -> stored-width mapping adds a +5/6 padding,
#  matching Excel's own column-width formula, so we back that offset out)
$ws.Columns.Item(4).ColumnWidth = 15 - 5/6
